$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in missing "+" marks in the grid
$ws.Range("D5").Value = "+"
$ws.Range("D7").Value = "+"
$ws.Range("F7").Value = "+"
$ws.Range("D9").Value = "+"
$ws.Range("F11").Value = "+"

# 2. Add new row 13 content
$ws.Range("A13").Value = "log.L.5.fl.19.conf.5.np.9728.dens.0.01`n  log.L.5.fl.20.conf.1.np.10240.dens.0.01`n  log.L.5.fl.20.conf.2.np.10240.dens.0.01`n  log.L.5.fl.20.conf.2.np.10240.dens.0.05`n  log.L.5.fl.21.conf.3.np.10752.dens.0.05"
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Range("A13").WrapText = $true
$ws.Range("A13").WrapText = $false

$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "``"

# 3. Update selection
$ws.Range("D10").Select()
